# Commit: "Sync attendance_reports, modules_schedules, and assets from main
# repo - 2026-01-31 04:00:58"
#
# The diff swaps the order of the two names recorded in the "Recorded By"
# column (column G) wherever the cell literally reads
# "dnasr281@gmail.com, System", turning it into "System, dnasr281@gmail.com".
# This happens on 78 rows of the "Session Analysis Results" sheet; every
# other cell (including other G-column values such as a lone
# "dnasr281@gmail.com") is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldText = "dnasr281@gmail.com, System"
$newText = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    if ($cell.Value() -eq $oldText) {
        $cell.Value = $newText
    }
}
